# TC02_Canine_Filter_Breed-AmerStaffd.xlsx
# "corrected ICDC Breed 1-14 scripts"
#
# The FilesTab (row 4) Cypher query dropped the `f.file_type` and
# `demo.breed` RETURN columns. Update the query text in B4, let the
# wrapped-text row shrink to its new (shorter) height, and leave the
# selection on B4 (mirroring where the cursor ended up after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['American Staffordshire Terrier']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newQuery

# The text shrank (two RETURN lines removed) so the wrapped row is shorter now.
$ws.Rows.Item(4).RowHeight = 217.5

# Cursor/selection ends up on B4 after the edit.
$ws.Range("B4").Select()
